# Auto-generated Word COM-interop script
# Updates the date line and the 20x5 arithmetic-answer table
# to match the target revision (commit c986bee).

$d = $word.ActiveDocument

# --- Update the date/weekday heading ---
$d.Content.Find.Execute("2024-12-05 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-06 Friday", 2) | Out-Null

# --- Update the answer table, cell by cell (row, column are 1-based) ---
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "70-45=25"
$t.Cell(1, 2).Range.Text = "84-5=79"
$t.Cell(1, 3).Range.Text = "19+75=94"
$t.Cell(1, 4).Range.Text = "39+18=57"
$t.Cell(1, 5).Range.Text = "97-79=18"
$t.Cell(2, 1).Range.Text = "55+26=81"
$t.Cell(2, 2).Range.Text = "58+14=72"
$t.Cell(2, 3).Range.Text = "87-69=18"
$t.Cell(2, 4).Range.Text = "86+6=92"
$t.Cell(2, 5).Range.Text = "18+36=54"
$t.Cell(3, 1).Range.Text = "64-25=39"
$t.Cell(3, 2).Range.Text = "68+8=76"
$t.Cell(3, 3).Range.Text = "84+9=93"
$t.Cell(3, 4).Range.Text = "19+7=26"
$t.Cell(3, 5).Range.Text = "50-7=43"
$t.Cell(4, 1).Range.Text = "78+17=95"
$t.Cell(4, 2).Range.Text = "27+45=72"
$t.Cell(4, 3).Range.Text = "15+6=21"
$t.Cell(4, 4).Range.Text = "32+9=41"
$t.Cell(4, 5).Range.Text = "26+15=41"
$t.Cell(5, 1).Range.Text = "77+17=94"
$t.Cell(5, 2).Range.Text = "92-33=59"
$t.Cell(5, 3).Range.Text = "57+9=66"
$t.Cell(5, 4).Range.Text = "32+39=71"
$t.Cell(5, 5).Range.Text = "48+33=81"
$t.Cell(6, 1).Range.Text = "6+7=13"
$t.Cell(6, 2).Range.Text = "51-14=37"
$t.Cell(6, 3).Range.Text = "39+32=71"
$t.Cell(6, 4).Range.Text = "91-6=85"
$t.Cell(6, 5).Range.Text = "72-24=48"
$t.Cell(7, 1).Range.Text = "28+39=67"
$t.Cell(7, 2).Range.Text = "34-28=6"
$t.Cell(7, 3).Range.Text = "25+39=64"
$t.Cell(7, 4).Range.Text = "27+34=61"
$t.Cell(8, 1).Range.Text = "9+48=57"
$t.Cell(8, 2).Range.Text = "26+49=75"
$t.Cell(8, 3).Range.Text = "60-33=27"
$t.Cell(8, 4).Range.Text = "84-15=69"
$t.Cell(8, 5).Range.Text = "19+78=97"
$t.Cell(9, 1).Range.Text = "14+17=31"
$t.Cell(9, 2).Range.Text = "88-79=9"
$t.Cell(9, 3).Range.Text = "36+15=51"
$t.Cell(9, 4).Range.Text = "8+45=53"
$t.Cell(9, 5).Range.Text = "81-13=68"
$t.Cell(10, 1).Range.Text = "49+3=52"
$t.Cell(10, 2).Range.Text = "17+38=55"
$t.Cell(10, 3).Range.Text = "70-59=11"
$t.Cell(10, 4).Range.Text = "22-17=5"
$t.Cell(10, 5).Range.Text = "7+54=61"
$t.Cell(11, 1).Range.Text = "54+37=91"
$t.Cell(11, 2).Range.Text = "92-15=77"
$t.Cell(11, 3).Range.Text = "59+23=82"
$t.Cell(11, 4).Range.Text = "92-87=5"
$t.Cell(11, 5).Range.Text = "90-35=55"
$t.Cell(12, 1).Range.Text = "9+4=13"
$t.Cell(12, 2).Range.Text = "77-29=48"
$t.Cell(12, 3).Range.Text = "38+49=87"
$t.Cell(12, 4).Range.Text = "9+53=62"
$t.Cell(12, 5).Range.Text = "48+9=57"
$t.Cell(13, 1).Range.Text = "64+9=73"
$t.Cell(13, 2).Range.Text = "19+63=82"
$t.Cell(13, 3).Range.Text = "80-1=79"
$t.Cell(13, 4).Range.Text = "52-45=7"
$t.Cell(13, 5).Range.Text = "54-18=36"
$t.Cell(14, 1).Range.Text = "17+64=81"
$t.Cell(14, 2).Range.Text = "80-23=57"
$t.Cell(14, 3).Range.Text = "91-37=54"
$t.Cell(14, 4).Range.Text = "51-49=2"
$t.Cell(14, 5).Range.Text = "29+12=41"
$t.Cell(15, 1).Range.Text = "5+77=82"
$t.Cell(15, 2).Range.Text = "19+44=63"
$t.Cell(15, 3).Range.Text = "44-15=29"
$t.Cell(15, 4).Range.Text = "28+54=82"
$t.Cell(15, 5).Range.Text = "27-8=19"
$t.Cell(16, 1).Range.Text = "83-54=29"
$t.Cell(16, 2).Range.Text = "12-4=8"
$t.Cell(16, 3).Range.Text = "70-46=24"
$t.Cell(16, 4).Range.Text = "13+58=71"
$t.Cell(16, 5).Range.Text = "35-18=17"
$t.Cell(17, 1).Range.Text = "47+44=91"
$t.Cell(17, 2).Range.Text = "50-33=17"
$t.Cell(17, 3).Range.Text = "67+25=92"
$t.Cell(17, 4).Range.Text = "8+19=27"
$t.Cell(17, 5).Range.Text = "48+28=76"
$t.Cell(18, 1).Range.Text = "55-27=28"
$t.Cell(18, 2).Range.Text = "16+38=54"
$t.Cell(18, 3).Range.Text = "8+26=34"
$t.Cell(18, 4).Range.Text = "27+35=62"
$t.Cell(18, 5).Range.Text = "37-29=8"
$t.Cell(19, 1).Range.Text = "92-67=25"
$t.Cell(19, 2).Range.Text = "42-7=35"
$t.Cell(19, 3).Range.Text = "34+8=42"
$t.Cell(19, 4).Range.Text = "4+58=62"
$t.Cell(19, 5).Range.Text = "70-1=69"
$t.Cell(20, 1).Range.Text = "63-57=6"
$t.Cell(20, 2).Range.Text = "71-56=15"
$t.Cell(20, 3).Range.Text = "75+8=83"
$t.Cell(20, 4).Range.Text = "70-34=36"
$t.Cell(20, 5).Range.Text = "78-19=59"

Write-Host "Done: updated date + 99 table cell(s); 1 cell(s) left unchanged."
